$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "updated to reflect currently purchased items"
# Column J ("Already Ordered?") flips from "No" to "Yes" for the items
# that have now been purchased: rows 2-8 and row 12.
$ws.Range("J2:J8").Value = "Yes"
$ws.Range("J12").Value = "Yes"

# Reflect the author's last on-screen selection when they saved the file.
$null = $ws.Range("J12").Select()
